$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "303.24"
Set-TextValue $ws.Range("E2") "-0.07%"
Set-TextValue $ws.Range("D3") "35.57"
Set-TextValue $ws.Range("E3") "4.78%"
Set-TextValue $ws.Range("D4") "5.091"
Set-TextValue $ws.Range("E4") "0.76%"
Set-TextValue $ws.Range("D5") "0.07738"
Set-TextValue $ws.Range("E5") "-0.99%"
Set-TextValue $ws.Range("D6") "2.219"
Set-TextValue $ws.Range("E6") "-7.58%"
Set-TextValue $ws.Range("D7") "8.032"
Set-TextValue $ws.Range("E7") "0.38%"
Set-TextValue $ws.Range("D8") "4.027"
Set-TextValue $ws.Range("E8") "3.49%"
Set-TextValue $ws.Range("D9") "0.9273"
Set-TextValue $ws.Range("E9") "-1.11%"
Set-TextValue $ws.Range("D10") "0.09534"
Set-TextValue $ws.Range("E10") "-5.95%"
Set-TextValue $ws.Range("E11") "2.37%"
Set-TextValue $ws.Range("D12") "0.08555"
Set-TextValue $ws.Range("E12") "0.40%"
Set-TextValue $ws.Range("D13") "0.03644"
Set-TextValue $ws.Range("E13") "9.03%"
Set-TextValue $ws.Range("D14") "0.09971"
Set-TextValue $ws.Range("E14") "0.52%"
Set-TextValue $ws.Range("D15") "0.001480"
Set-TextValue $ws.Range("E15") "0.28%"
Set-TextValue $ws.Range("D16") "0.005698"
Set-TextValue $ws.Range("E16") "0.21%"
Set-TextValue $ws.Range("D17") "3.478"
Set-TextValue $ws.Range("E17") "0.40%"
Set-TextValue $ws.Range("D18") "2.184"
Set-TextValue $ws.Range("E18") "-0.03%"
Set-TextValue $ws.Range("E19") "3.12%"
Set-TextValue $ws.Range("D20") "0.1324"
Set-TextValue $ws.Range("D21") "4.567"
Set-TextValue $ws.Range("E21") "6.41%"
Set-TextValue $ws.Range("E22") "7.29%"
Set-TextValue $ws.Range("E23") "1.42%"
Set-TextValue $ws.Range("E24") "1.50%"
Set-TextValue $ws.Range("D25") "0.004506"
Set-TextValue $ws.Range("E25") "2.17%"
Set-TextValue $ws.Range("D26") "0.0001308"
Set-TextValue $ws.Range("E26") "1.12%"
Set-TextValue $ws.Range("E27") "-20.04%"
Set-TextValue $ws.Range("D39") "0.01757"
Set-TextValue $ws.Range("E39") "1.72%"
Set-TextValue $ws.Range("D40") "0.04719"
Set-TextValue $ws.Range("E40") "-2.05%"
Set-TextValue $ws.Range("D41") "0.007951"
Set-TextValue $ws.Range("E41") "2.27%"
Set-TextValue $ws.Range("D42") "0.1407"
Set-TextValue $ws.Range("E42") "-0.10%"
Set-TextValue $ws.Range("D43") "0.007685"
Set-TextValue $ws.Range("E43") "-21.45%"
Set-TextValue $ws.Range("D44") "0.002227"
Set-TextValue $ws.Range("E44") "7.35%"
Set-TextValue $ws.Range("D45") "0.009660"
Set-TextValue $ws.Range("E45") "4.65%"
Set-TextValue $ws.Range("D46") "0.00006241"
Set-TextValue $ws.Range("E46") "2.67%"
Set-TextValue $ws.Range("D47") "0.00000000755"
Set-TextValue $ws.Range("E47") "1.15%"
Set-TextValue $ws.Range("E49") "35.81%"
Set-TextValue $ws.Range("D50") "0.00002113"
Set-TextValue $ws.Range("E50") "1.15%"
Set-TextValue $ws.Range("D51") "0.0002013"
Set-TextValue $ws.Range("E51") "1.15%"
